# Use case "Visualizzazione modello": fix grammatical/logic error.
# The sentence previously read "... segnala che il server non è stato
# possibile effettuare la comunicazione di sistema." The reference to
# "il server" is removed so it reads "... segnala che non è stato
# possibile effettuare la comunicazione di sistema."

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "il server non è ",  # FindText
    $true,                # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                     # Wrap (wdFindContinue)
    $false,               # Format
    "non è ",             # ReplaceWith
    2                      # Replace (wdReplaceAll)
)
